# Update the test-data username value and move the active selection,
# matching the authoring change reflected in the workbook diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A2 held the shared string "qatitans31" -> bump to "qatitans40"
$ws.Range("A2").Value = "qatitans40"

# Active selection moves from C2 to A2
$null = $ws.Range("A2").Select()
